$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("L3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("M3").Value = 2021
$ws.Range("N3").Value = 2022

$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("L4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("M4").Value = 6.18
$ws.Range("N4").Value = 6.18
